# Fixed 3 watchlist test cases
#
# 1. "IAM010" sheet (CHARACTER LENGTH watchlist): bump the three character-length
#    seed values (61/62/63 -> 246/247/248 and 11 -> 111 twice) and update the
#    truncated-error-message text to reflect the new 255-character limit.
# 2. Active tab moves from "Test Cases" to "IAM010", and the selection on each of
#    those two sheets changes accordingly.

$wb = $excel.ActiveWorkbook

# --- IAM010: update the watchlist test-case values -------------------------
$ws = $wb.Worksheets.Item("IAM010")
$ws.Range("A2").Value = 246
$ws.Range("A3").Value = 247
$ws.Range("A4").Value = 248
$ws.Range("C4").Value = "Please enter no more than 255 characters."
$ws.Range("A5").Value = 111
$ws.Range("A6").Value = 111

# --- Move the active tab from "Test Cases" to "IAM010" ---------------------
$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Activate()
$wsTestCases.Range("A11").Select()

$ws.Activate()
$ws.Range("A8").Select()
